$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the existing 10 data rows (rows 2-11) into rows 12-21,
# keeping year/employee_id values, but only keeping the accident_type
# text for the 2017 rows (and for employee 5's 2016 row).
$data = @(
    @(2016, 1, $null),
    @(2017, 1, "Mild"),
    @(2016, 2, $null),
    @(2017, 2, "Mild"),
    @(2016, 4, $null),
    @(2017, 4, "Mild"),
    @(2016, 5, "Moderate"),
    @(2017, 5, "Severe"),
    @(2016, 7, $null),
    @(2017, 7, "Moderate")
)

$row = 12
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    if ($entry[2]) {
        $ws.Cells.Item($row, 3).Value = $entry[2]
    }
    $row = $row + 1
}

$ws.Range("E14").Select()
